$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @("08/27","11/16","02/17","07/17","11/24","03/08","01/28","04/20","05/05","01/09","01/30","08/25","09/09","10/28","04/30","11/07")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("C$row").Value = $values[$i]
}

$ws.Range("C2").Select()
